$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cell values (content reshuffle per source data update) ---
$ws.Range("B10").Value = "8711623 - Denize Kalempa"
$ws.Range("C10").Value = "8711623 - Denize Kalempa"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Statics and dynamics of fluids, oscillations and mechanical waves, ideal gas,temperature, heat and the laws of thermodynamics."
$ws.Range("C14").Value = "Statics and dynamics of fluids, oscillations and mechanical waves, ideal gas,temperature, heat and the laws of thermodynamics."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1) Fluids at rest: pressure, Pascal’s principle , Archimedes’ principle, surfacetension and capillarity;2) Fluids in motion: flow rate, ideal fluids, the equation of continuity, Bernoulli’s equation, viscosity and the Hagen-Poiseuille law;3) Oscillation: simple harmonic motion, damped and forced oscillations, resonance;4) Waves: transverse and longitudinal, wave equation, superposition, interference, standing waves, sound waves, intensity and sound level, beats, Doppler effect;5) Temperature and heat: definitions, zeroth Law of thermodynamics, thermal expansion, absorption of heat by solids and liquids, heat and work, heat transfer mechanisms, ideal gases, specific heat and degrees of freedom for an ideal gas;6) Thermodynamics: the first law of thermodynamics, reversible and irreversible processes, heat engines and efficiency, entropy, the second law of thermodynamics."
$ws.Range("C16").Value = "1) Fluids at rest: pressure, Pascal’s principle , Archimedes’ principle, surfacetension and capillarity;2) Fluids in motion: flow rate, ideal fluids, the equation of continuity, Bernoulli’s equation, viscosity and the Hagen-Poiseuille law;3) Oscillation: simple harmonic motion, damped and forced oscillations, resonance;4) Waves: transverse and longitudinal, wave equation, superposition, interference, standing waves, sound waves, intensity and sound level, beats, Doppler effect;5) Temperature and heat: definitions, zeroth Law of thermodynamics, thermal expansion, absorption of heat by solids and liquids, heat and work, heat transfer mechanisms, ideal gases, specific heat and degrees of freedom for an ideal gas;6) Thermodynamics: the first law of thermodynamics, reversible and irreversible processes, heat engines and efficiency, entropy, the second law of thermodynamics."
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8711623 - Denize Kalempa"
$ws.Range("C18").Value = "8711623 - Denize Kalempa"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("B24").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1018 -  Física I  (Requisito fraco)`n"

# --- Clear cells that become empty (without shifting neighbours) ---
$ws.Range("B17:C17").ClearContents()
$ws.Range("B22:C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Row height adjustments ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# --- Remove the now-empty trailing row 25 ---
$ws.Rows.Item(25).Delete()
